# Commit: "spr 20240227 add ToneBurst"
#
# A new configuration row (ID 102.1) is inserted after the existing ID=102
# ("TB-ThreeScales_Variance") row. It is a ToneBurst-style variant of that
# same condition: same Info/folderName/ParentFolderName etc., but a newer
# Version_Date and a new sigma ("S") value list. All rows from the old
# row 5 onward shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a blank row at position 5 (row 2 is the type-info row, rows 3-10
# are data rows 101,102,103,104,<blank>,201,202,203). This shifts the old
# row 5 (ID 103) and everything below it down by one, and the new blank
# row 5 inherits the formatting of what is now row 6 (old row 5), which is
# exactly the formatting/border pattern the target row needs.
$ws.Rows.Item(5).Insert()

# Populate the new row 5 as a near-duplicate of row 4 (ID 102,
# TB-ThreeScales_Variance), but with ID 102.1, an updated Version_Date, and
# a new sigma list (4,4;4,8;4,16;4,32;4,64) reflecting the new ToneBurst
# variant.
$ws.Cells.Item(5, 1).Value = 102.1                              # A: ID
$ws.Cells.Item(5, 2).Value = "2024/2/21"                          # B: Version_Date
$ws.Cells.Item(5, 3).Value = 97656                                # C: fs
$ws.Cells.Item(5, 4).Value = "TB-ThreeScales_Variance"            # D: Info
$ws.Cells.Item(5, 5).Value = 1                                    # E: saveMat
$ws.Cells.Item(5, 6).Value = "MNP_TB-ThreeScales_Variance"        # F: folderName
$ws.Cells.Item(5, 7).Value = "MonkeyNeuroPixels"                  # G: ParentFolderName
$ws.Cells.Item(5, 8).Value = 0.2                                  # H: clickDur
$ws.Cells.Item(5, 9).Value = 0.5                                  # I: Amp
$ws.Cells.Item(5, 13).Value = "pulse"                             # M: clickType
$ws.Cells.Item(5, 15).Value = 2000                                # O: S1Dur
$ws.Cells.Item(5, 16).Value = "2000"                              # P: S2Dur
$ws.Cells.Item(5, 17).Value = "4"                                 # Q: ICIBase
$ws.Cells.Item(5, 18).Value = "2"                                 # R: ratio
$ws.Cells.Item(5, 19).Value = "4,4;4,8;4,16;4,32;4,64"            # S: sigma
$ws.Cells.Item(5, 20).Value = "0.25,2"                            # T: ICIRangeRatio
$ws.Cells.Item(5, 21).Value = "0"                                 # U: skewBase
$ws.Cells.Item(5, 22).Value = "0"                                 # V: skewType
$ws.Cells.Item(5, 23).Value = "0"                                 # W: repNs
$ws.Cells.Item(5, 26).Value = "1"                                 # Z: lastClick
$ws.Cells.Item(5, 27).Value = 0                                   # AA: repHead
$ws.Cells.Item(5, 28).Value = 0                                   # AB: repTail
$ws.Cells.Item(5, 29).Value = "0.9,1.1"                           # AC: repRatio
$ws.Cells.Item(5, 40).Value = "Irreg"                             # AN: soundType
$ws.Cells.Item(5, 42).Value = "@ThreeScales_Variance_Gen"         # AP: GenFcn

# The "insert row" operation above inherited formatting from the row that
# is now 6 rows below (old row 5), which lacks the left/right borders that
# this sheet uses on columns O and AC to visually separate sections. Match
# the canonical border pattern used by every other data row.
$ws.Cells.Item(5, 15).Borders.Item(7).LineStyle = 1    # O5: xlEdgeLeft
$ws.Cells.Item(5, 29).Borders.Item(10).LineStyle = 1   # AC5: xlEdgeRight

# The row that used to be row 6 (ID 104, "TB-ThreeScales_Skew_8ms") is now
# row 7 after the insert; its ID is renumbered to 103.1 to reflect that it
# is a sibling/variant of the row above it (ID 103) rather than a
# standalone ID 104.
$ws.Cells.Item(7, 1).Value = 103.1                                # A7: ID
